$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): force text so numeric-looking strings (e.g. "2.40")
# keep their literal formatting instead of being coerced to a Double.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "65.308.10"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.397.87"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "592.48"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "142.02"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.402.08"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.466"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.134"
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "7.86"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.404"
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "3.989.37"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "29.66"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.398.70"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "65.454.53"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "10.33"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.08"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.73"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "416.22"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.579"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "77.02"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.542.36"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0000109"
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.17"
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.79"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.40"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.161"
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "24.53"
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.70"
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.51"
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "7.53"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "172.83"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0861"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.03"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.866"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "45.46"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "26.50"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.16"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.06"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.26"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.915"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.232"
$c.Style = "Normal"

# Volume(1h) column (E): plain text assignment, values already
# contain "%" / leading-minus so Excel stores them as text natively.
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("E6").Value = "  -3.55%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("E9").Value = "  -3.43%  "
$ws.Range("E10").Value = "  -5.25%  "
$ws.Range("E11").Value = "  +5.45%  "
$ws.Range("E12").Value = "  -4.17%  "
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("E14").Value = "  -5.88%  "
$ws.Range("E15").Value = "  -6.23%  "
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("E19").Value = "  +3.69%  "
$ws.Range("E20").Value = "  -5.10%  "
$ws.Range("E21").Value = "  -3.48%  "
$ws.Range("E22").Value = "  -5.05%  "
$ws.Range("E23").Value = "  -4.66%  "
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -1.48%  "
$ws.Range("E27").Value = "  -8.74%  "
$ws.Range("E28").Value = "  -6.86%  "
$ws.Range("E29").Value = "  -7.46%  "
$ws.Range("E30").Value = "  -2.80%  "
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("E33").Value = "  -8.74%  "
$ws.Range("E34").Value = "  -3.15%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -5.49%  "
$ws.Range("E37").Value = "  -9.54%  "
$ws.Range("E38").Value = "  -4.68%  "
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  -2.69%  "
$ws.Range("E42").Value = "  -6.21%  "
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("E44").Value = "  -11.77%  "
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("E46").Value = "  -8.19%  "
$ws.Range("E47").Value = "  -5.40%  "
$ws.Range("E48").Value = "  -5.43%  "
$ws.Range("E49").Value = "  -7.61%  "
$ws.Range("E50").Value = "  -7.14%  "
$ws.Range("E51").Value = "  -5.11%  "
